# Update column F ("dSF") values for a set of rows.
# These reflect a repull/recalculation of the data (per commit message:
# "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -2
    8  = -2
    11 = 0
    13 = 2
    18 = 2
    21 = 0
    28 = -1
    29 = 1
    34 = 2
    35 = 4
    37 = 5
    38 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
